$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F24").Value = 13
$ws.Range("G24").Value = 458.77
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 1259.52
$ws.Range("B32").Value = 12689.37
$ws.Range("F70").Value = 15
$ws.Range("G70").Value = 2024.25
$ws.Range("F71").Value = 322
$ws.Range("G71").Value = 20511.4
$ws.Range("F77").Value = 254
$ws.Range("G77").Value = 11871.96
$ws.Range("F80").Value = 8
$ws.Range("G80").Value = 1968.56
$ws.Range("F83").Value = 115
$ws.Range("G83").Value = 17327.05
$ws.Range("F85").Value = 141
$ws.Range("G85").Value = 19002.57
$ws.Range("B90").Value = 179507.3
$ws.Range("F95").Value = 3
$ws.Range("G95").Value = 424.29
$ws.Range("B96").Value = 424.29
$ws.Range("B112").Value = 64350
$ws.Range("E112").Value = 70.63
$ws.Range("F112").Value = 2
$ws.Range("G112").Value = 132.88
$ws.Range("B113").Value = 57756
$ws.Range("E113").Value = 79.37
$ws.Range("F113").Value = -100
$ws.Range("G113").Value = -6644
$ws.Range("F115").Value = 202
$ws.Range("G115").Value = 19555.62
$ws.Range("B117").Value = 13413.5
$ws.Range("F144").Value = 1039
$ws.Range("G144").Value = 8779.549999999999
$ws.Range("B147").Value = 14403.01
$ws.Range("F149").Value = 228
$ws.Range("G149").Value = 14774.4
$ws.Range("F150").Value = 35
$ws.Range("G150").Value = 1627.15
$ws.Range("B156").Value = 31454.5
$ws.Range("B192").Value = 64973
$ws.Range("E192").Value = 35.4
$ws.Range("F192").Value = 2
$ws.Range("G192").Value = 66.59999999999999
$ws.Range("B193").Value = 48706
$ws.Range("E193").Value = 39.8
$ws.Range("F193").Value = -144
$ws.Range("G193").Value = -4795.2
$ws.Range("F199").Value = 23
$ws.Range("G199").Value = 5689.74
$ws.Range("F203").Value = 58
$ws.Range("G203").Value = 1169.28
$ws.Range("B216").Value = 39768.05
$ws.Range("B219").Value = 63565
$ws.Range("E219").Value = 109.19
$ws.Range("F219").Value = 60
$ws.Range("G219").Value = 6162.6
$ws.Range("B220").Value = 61610
$ws.Range("E220").Value = 122.71
$ws.Range("F220").Value = -58
$ws.Range("G220").Value = -5957.18
$ws.Range("B227").Value = 63520
$ws.Range("E227").Value = 153.4
$ws.Range("F227").Value = 66
$ws.Range("G227").Value = 9522.48
$ws.Range("B228").Value = 55373
$ws.Range("E228").Value = 163.62
$ws.Range("F228").Value = -94
$ws.Range("G228").Value = -13562.32
$ws.Range("F229").Value = 61
$ws.Range("G229").Value = 8752.280000000001
$ws.Range("F234").Value = 40
$ws.Range("G234").Value = 2052.8
$ws.Range("B243").Value = 63560
$ws.Range("E243").Value = 134.87
$ws.Range("F243").Value = 1
$ws.Range("G243").Value = 126.86
$ws.Range("B244").Value = 60325
$ws.Range("E244").Value = 151.57
$ws.Range("F244").Value = -102
$ws.Range("G244").Value = -12939.72
$ws.Range("F247").Value = 141
$ws.Range("G247").Value = 14651.31
$ws.Range("F249").Value = 138
$ws.Range("G249").Value = 19019.16
$ws.Range("F250").Value = 9
$ws.Range("G250").Value = 4449.42
$ws.Range("F251").Value = 0
$ws.Range("G251").Value = 0
$ws.Range("F252").Value = 1
$ws.Range("G252").Value = 21.03
$ws.Range("F255").Value = 554
$ws.Range("G255").Value = 94916.82000000001
$ws.Range("B260").Value = 187519.3
$ws.Range("F280").Value = 134
$ws.Range("G280").Value = 22664.76
$ws.Range("F283").Value = 40
$ws.Range("G283").Value = 13658.8
$ws.Range("F302").Value = 44
$ws.Range("G302").Value = 9279.16
$ws.Range("B304").Value = 174236.87
$ws.Range("F320").Value = 45
$ws.Range("G320").Value = 3089.25
$ws.Range("B330").Value = 27564.18
$ws.Range("F334").Value = 192
$ws.Range("G334").Value = 9949.440000000001
$ws.Range("F339").Value = 0
$ws.Range("G339").Value = 0
$ws.Range("F345").Value = 52
$ws.Range("G345").Value = 3193.32
$ws.Range("B346").Value = 25832.85
$ws.Range("B442").Value = 53319
$ws.Range("E442").Value = 310.64
$ws.Range("F442").Value = -6
$ws.Range("G442").Value = -1643.52
$ws.Range("B443").Value = 64810
$ws.Range("E443").Value = 291.22
$ws.Range("F443").Value = 4
$ws.Range("G443").Value = 1095.68
$ws.Range("F453").Value = 21
$ws.Range("G453").Value = 556.71
$ws.Range("F454").Value = 49
$ws.Range("G454").Value = 1673.35
$ws.Range("B460").Value = 13609.42
$ws.Range("B463").Value = 64833
$ws.Range("E463").Value = 34.9
$ws.Range("F463").Value = 95
$ws.Range("G463").Value = 3118.85
$ws.Range("B464").Value = 60025
$ws.Range("E464").Value = 37.22
$ws.Range("F464").Value = -98
$ws.Range("G464").Value = -3217.34
$ws.Range("F486").Value = 75
$ws.Range("G486").Value = 6625.5
$ws.Range("F487").Value = 0
$ws.Range("G487").Value = 0
$ws.Range("B488").Value = 30455.58
$ws.Range("F508").Value = 56
$ws.Range("G508").Value = 5820.64
$ws.Range("F509").Value = 215
$ws.Range("G509").Value = 17281.7
$ws.Range("B510").Value = 23102.34
$ws.Range("F550").Value = 0
$ws.Range("G550").Value = 0
$ws.Range("F555").Value = 19
$ws.Range("G555").Value = 1321.64
$ws.Range("B560").Value = 4694.38
$ws.Range("F577").Value = 59
$ws.Range("G577").Value = 2536.41
$ws.Range("F580").Value = 56
$ws.Range("G580").Value = 3191.44
$ws.Range("F582").Value = 33
$ws.Range("G582").Value = 1880.67
$ws.Range("B583").Value = 16333.02
$ws.Range("F599").Value = 1620
$ws.Range("G599").Value = 264238.2
$ws.Range("F601").Value = 408
$ws.Range("G601").Value = 115410.96
$ws.Range("B606").Value = 428231.71
$ws.Range("F612").Value = 32
$ws.Range("G612").Value = 1311.68
$ws.Range("B618").Value = 43689.86
$ws.Range("B619").Value = 1738534.67
$ws.Range("B620").Value = 1738534.67

Write-Output "Applied 164 cell updates"
